# Add a new "Swiss" worksheet with Switzerland market test data,
# based on a copy of the existing "Czech" worksheet.

$wb = $excel.ActiveWorkbook

# Duplicate the "Czech" sheet and place the copy immediately after it.
$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)

# The newly created copy becomes the active sheet and is named "Czech (2)"
# by default - rename it to "Swiss".
$swiss = $wb.Worksheets.Item("Czech (2)")
$swiss.Name = "Swiss"

# Update the market-specific values for Switzerland.
# Order matters: new shared strings are appended in the order they are
# first used, so set them in the same order as in the target workbook.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("A16").Value = "P32AR-CH"
$swiss.Range("A17").Value = "P32DR-CH"
$swiss.Range("B4").Value = "NGC-3476/T2645"
